$wb = $excel.ActiveWorkbook

# ----- Sheet "Size" -----
$wsSize = $wb.Worksheets.Item("Size")
$wsSize.Range("G2").Value = 320.37394354013549
$wsSize.Range("G3").Value = 1620.946024170172
$wsSize.Range("G4").Value = 18.714822794742179

$wsSize.Range("C5").Value = 244.98844740528031
$wsSize.Range("D5").Value = 1069.252757002477
$wsSize.Range("E5").Value = 54.722222994710577
$wsSize.Range("F5").Value = 47.442304286471263
$wsSize.Range("G5").Value = 1416.4057316889391

# ----- Sheet "Cost" -----
$wsCost = $wb.Worksheets.Item("Cost")
$wsCost.Range("I2").Value = 9.9840096059080814

$wsCost.Range("I3").Value = 0.64074788708027108
$wsCost.Range("I4").Value = 0.89152031329359449
$wsCost.Range("I5").Value = 0.0037429645589484349

$wsCost.Range("E6").Value = 0.02449884474052803
$wsCost.Range("F6").Value = 0.10692527570024769
$wsCost.Range("G6").Value = 0.0054722222994710579
$wsCost.Range("H6").Value = 0.0047442304286471251
$wsCost.Range("I6").Value = 0.14164057316889389

$wsCost.Range("I7").Value = 0.1788789342984472
$wsCost.Range("I8").Value = 0.24888759957377349
$wsCost.Range("I9").Value = 0.0052246564126217391

$wsCost.Range("E10").Value = 0.0051295454816456159
$wsCost.Range("F10").Value = 0.022387915456869651
$wsCost.Range("G10").Value = 0.0011457688502501671
$wsCost.Range("H10").Value = 0.00099334258479928596
$wsCost.Range("I10").Value = 0.02965657237356472

$wsCost.Range("I11").Value = 0.61993146406142541

$wsCost.Range("I12").Value = 0.13939342809793881

$wsCost.Range("E13").Value = 1.9334535389604219
$wsCost.Range("F13").Value = 5.6675463408115663
$wsCost.Range("G13").Value = 0.082659900189235955
$wsCost.Range("H13").Value = 0.020656897088802501
$wsCost.Range("I13").Value = 7.7043166770500262

# ----- Sheet "Indicators" -----
$wsInd = $wb.Worksheets.Item("Indicators")
$wsInd.Range("C2").Value = 960.46647168614982
$wsInd.Range("D2").Value = 2787.02527677542
$wsInd.Range("E2").Value = 3747.4917484615698

$wsInd.Range("E3").Value = 0.76992550822477179

$wsInd.Range("E4").Value = 0.23007449177522829

$wsInd.Range("C5").Value = 0.82901061634023565
$wsInd.Range("D5").Value = 0.87878523447900447
$wsInd.Range("E5").Value = 0.86602820789371726

$wsInd.Range("E6").Value = 0.22038916861858751
